# Fixed a bug in calcWaysWinsInReels2
# The data rows (2-25) containing reel-strip symbol/weight data got
# reordered. Re-arrange the rows so that each row's A:F data lands on its
# corrected row, matching the fixed ordering produced by
# calcWaysWinsInReels2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot current values of rows 2..25, columns A..F, before overwriting
# anything (rows later in the order may depend on rows earlier in the
# sheet, so we must capture everything up-front).
$snapshot = @{}
for ($r = 2; $r -le 25; $r++) {
    $rowVals = @{}
    foreach ($col in @("A","B","C","D","E","F")) {
        $rowVals[$col] = $ws.Range("$col$r").Value()
    }
    $snapshot[$r] = $rowVals
}

# Mapping of destination row -> source row (captured in the snapshot above)
$mapping = @{
    2  = 14
    3  = 8
    4  = 7
    5  = 15
    6  = 10
    7  = 4
    8  = 5
    9  = 6
    10 = 2
    11 = 13
    12 = 11
    13 = 9
    14 = 3
    15 = 12
    16 = 20
    17 = 18
    18 = 16
    19 = 21
    20 = 19
    21 = 17
    22 = 23
    23 = 22
    24 = 24
    25 = 25
}

foreach ($dest in $mapping.Keys) {
    $src = $mapping[$dest]
    $srcVals = $snapshot[$src]
    foreach ($col in @("A","B","C","D","E","F")) {
        $ws.Range("$col$dest").Value = $srcVals[$col]
    }
}
